$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(112, 8).Value = 4193.4517
$ws.Cells.Item(112, 10).Value = 4389.552
$ws.Cells.Item(112, 12).Value = 13168.656
$ws.Cells.Item(112, 14).Value = -15384.656
$ws.Cells.Item(137, 8).Value = 359592.94
$ws.Cells.Item(137, 9).Value = 386861.2
$ws.Cells.Item(137, 10).Value = 5105.5
$ws.Cells.Item(137, 11).Value = 1160583.6
$ws.Cells.Item(137, 12).Value = 15316.5
$ws.Cells.Item(137, 13).Value = -1158033.6
$ws.Cells.Item(137, 14).Value = -20416.5
$ws.Cells.Item(138, 8).Value = 5831.761
$ws.Cells.Item(138, 9).Value = 1542.4783
$ws.Cells.Item(138, 10).Value = 8073.886
$ws.Cells.Item(138, 11).Value = 4627.4349
$ws.Cells.Item(138, 12).Value = 24221.658
$ws.Cells.Item(138, 13).Value = 512.5650999999998
$ws.Cells.Item(138, 14).Value = -34501.658

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14873.324
$ws.Cells.Item(32, 9).Value = 15268.597
$ws.Cells.Item(32, 10).Value = 13548
$ws.Cells.Item(32, 11).Value = 15268.597
$ws.Cells.Item(32, 12).Value = 13548
$ws.Cells.Item(32, 13).Value = -14981.597
$ws.Cells.Item(32, 14).Value = -14122
$ws.Cells.Item(88, 8).Value = 47022.816
$ws.Cells.Item(88, 9).Value = 337.5
$ws.Cells.Item(88, 10).Value = 73700.14
$ws.Cells.Item(88, 11).Value = 337.5
$ws.Cells.Item(88, 12).Value = 73700.14
$ws.Cells.Item(88, 13).Value = 68.5
$ws.Cells.Item(88, 14).Value = -74512.14
$ws.Cells.Item(91, 8).Value = 47022.816
$ws.Cells.Item(91, 9).Value = 337.5
$ws.Cells.Item(91, 10).Value = 73700.14
$ws.Cells.Item(91, 11).Value = 337.5
$ws.Cells.Item(91, 12).Value = 73700.14
$ws.Cells.Item(91, 13).Value = 1066.5
$ws.Cells.Item(91, 14).Value = -76508.14
$ws.Cells.Item(97, 8).Value = 1104.762
$ws.Cells.Item(97, 9).Value = 1125.8422
$ws.Cells.Item(97, 10).Value = 904.5
$ws.Cells.Item(97, 11).Value = 1125.8422
$ws.Cells.Item(97, 12).Value = 904.5
$ws.Cells.Item(97, 13).Value = -629.8422
$ws.Cells.Item(97, 14).Value = -1896.5
$ws.Cells.Item(102, 8).Value = 1097.2632
$ws.Cells.Item(102, 9).Value = 928.1875
$ws.Cells.Item(102, 10).Value = 1999
$ws.Cells.Item(102, 11).Value = 928.1875
$ws.Cells.Item(102, 12).Value = 1999
$ws.Cells.Item(102, 13).Value = 693.8125
$ws.Cells.Item(102, 14).Value = -5243
$ws.Cells.Item(122, 8).Value = 3860.375
$ws.Cells.Item(122, 9).Value = 3403.1277
$ws.Cells.Item(122, 10).Value = 5124.5293
$ws.Cells.Item(122, 11).Value = 10209.3831
$ws.Cells.Item(122, 12).Value = 15373.5879
$ws.Cells.Item(122, 13).Value = -7759.383099999999
$ws.Cells.Item(122, 14).Value = -20273.5879
$ws.Cells.Item(132, 8).Value = 13805.962
$ws.Cells.Item(132, 9).Value = 16351.795
$ws.Cells.Item(132, 10).Value = 6168.4614
$ws.Cells.Item(132, 11).Value = 49055.385
$ws.Cells.Item(132, 12).Value = 18505.3842
$ws.Cells.Item(132, 13).Value = -46525.385
$ws.Cells.Item(132, 14).Value = -23565.3842

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2005.2
$ws.Cells.Item(20, 9).Value = 1986.4286
$ws.Cells.Item(20, 10).Value = 2049
$ws.Cells.Item(20, 11).Value = 1986.4286
$ws.Cells.Item(20, 12).Value = 2049
$ws.Cells.Item(20, 13).Value = -1739.4286
$ws.Cells.Item(20, 14).Value = -2543
$ws.Cells.Item(42, 8).Value = 279999
$ws.Cells.Item(42, 10).Value = 279999
$ws.Cells.Item(42, 12).Value = 279999
$ws.Cells.Item(42, 14).Value = -280655
$ws.Cells.Item(86, 8).Value = 1151.2084
$ws.Cells.Item(86, 9).Value = 1004.94116
$ws.Cells.Item(86, 10).Value = 1506.4286
$ws.Cells.Item(86, 11).Value = 1004.94116
$ws.Cells.Item(86, 12).Value = 1506.4286
$ws.Cells.Item(86, 13).Value = 118.05884
$ws.Cells.Item(86, 14).Value = -3752.4286
$ws.Cells.Item(89, 8).Value = 1151.2084
$ws.Cells.Item(89, 9).Value = 1004.94116
$ws.Cells.Item(89, 10).Value = 1506.4286
$ws.Cells.Item(89, 11).Value = 5024.7058
$ws.Cells.Item(89, 12).Value = 7532.143
$ws.Cells.Item(89, 13).Value = 591.2942000000003
$ws.Cells.Item(89, 14).Value = -18764.143
$ws.Cells.Item(94, 8).Value = 943.2727
$ws.Cells.Item(94, 9).Value = 798.4074000000001
$ws.Cells.Item(94, 10).Value = 1595.1666
$ws.Cells.Item(94, 11).Value = 798.4074000000001
$ws.Cells.Item(94, 12).Value = 1595.1666
$ws.Cells.Item(94, 13).Value = -347.4074000000001
$ws.Cells.Item(94, 14).Value = -2497.1666
$ws.Cells.Item(134, 8).Value = 1731.3334
$ws.Cells.Item(134, 9).Value = 1451.5135
$ws.Cells.Item(134, 11).Value = 4354.5405
$ws.Cells.Item(134, 13).Value = -1819.5405

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2745.7637
$ws.Cells.Item(31, 9).Value = 1457.1818
$ws.Cells.Item(31, 10).Value = 7900.091
$ws.Cells.Item(31, 11).Value = 1457.1818
$ws.Cells.Item(31, 12).Value = 7900.091
$ws.Cells.Item(31, 13).Value = -1162.1818
$ws.Cells.Item(31, 14).Value = -8490.091
$ws.Cells.Item(34, 8).Value = 2745.7637
$ws.Cells.Item(34, 9).Value = 1457.1818
$ws.Cells.Item(34, 10).Value = 7900.091
$ws.Cells.Item(34, 11).Value = 1457.1818
$ws.Cells.Item(34, 12).Value = 7900.091
$ws.Cells.Item(34, 13).Value = -1255.1818
$ws.Cells.Item(34, 14).Value = -8304.091
$ws.Cells.Item(122, 8).Value = 4046.5454
$ws.Cells.Item(122, 9).Value = 2468.2856
$ws.Cells.Item(122, 11).Value = 7404.8568
$ws.Cells.Item(122, 13).Value = -4954.8568
$ws.Cells.Item(132, 8).Value = 30310530
$ws.Cells.Item(132, 9).Value = 35091020
$ws.Cells.Item(132, 10).Value = 34099.668
$ws.Cells.Item(132, 11).Value = 105273060
$ws.Cells.Item(132, 12).Value = 102299.004
$ws.Cells.Item(132, 13).Value = -105270530
$ws.Cells.Item(132, 14).Value = -107359.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(103, 8).Value = 622.1539
$ws.Cells.Item(103, 9).Value = 312.5
$ws.Cells.Item(103, 11).Value = 937.5
$ws.Cells.Item(103, 13).Value = -58.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 82621.28999999999
$ws.Cells.Item(80, 9).Value = 128324.875
$ws.Cells.Item(80, 10).Value = 21683.166
$ws.Cells.Item(80, 11).Value = 128324.875
$ws.Cells.Item(80, 12).Value = 21683.166
$ws.Cells.Item(80, 13).Value = -127326.875
$ws.Cells.Item(80, 14).Value = -23679.166
$ws.Cells.Item(83, 8).Value = 82621.28999999999
$ws.Cells.Item(83, 9).Value = 128324.875
$ws.Cells.Item(83, 10).Value = 21683.166
$ws.Cells.Item(83, 11).Value = 641624.375
$ws.Cells.Item(83, 12).Value = 108415.83
$ws.Cells.Item(83, 13).Value = -636632.375
$ws.Cells.Item(83, 14).Value = -118399.83
$ws.Cells.Item(102, 8).Value = 13520359
$ws.Cells.Item(102, 9).Value = 22735522
$ws.Cells.Item(102, 10).Value = 4785.8
$ws.Cells.Item(102, 11).Value = 22735522
$ws.Cells.Item(102, 12).Value = 4785.8
$ws.Cells.Item(102, 13).Value = -22733900
$ws.Cells.Item(102, 14).Value = -8029.8
$ws.Cells.Item(126, 8).Value = 3529.25
$ws.Cells.Item(126, 9).Value = 2038.1
$ws.Cells.Item(126, 11).Value = 6114.299999999999
$ws.Cells.Item(126, 13).Value = -3644.299999999999
$ws.Cells.Item(132, 8).Value = 43954.18
$ws.Cells.Item(132, 9).Value = 53635.25
$ws.Cells.Item(132, 10).Value = 5229.9
$ws.Cells.Item(132, 11).Value = 160905.75
$ws.Cells.Item(132, 12).Value = 15689.7
$ws.Cells.Item(132, 13).Value = -158375.75
$ws.Cells.Item(132, 14).Value = -20749.7

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(61, 8).Value = 1678.1428
$ws.Cells.Item(61, 9).Value = 1576.5385
$ws.Cells.Item(61, 10).Value = 2999
$ws.Cells.Item(61, 11).Value = 1576.5385
$ws.Cells.Item(61, 12).Value = 2999
$ws.Cells.Item(61, 13).Value = -1374.5385
$ws.Cells.Item(61, 14).Value = -3403
$ws.Cells.Item(113, 8).Value = 1678.1428
$ws.Cells.Item(113, 9).Value = 1576.5385
$ws.Cells.Item(113, 10).Value = 2999
$ws.Cells.Item(113, 11).Value = 1576.5385
$ws.Cells.Item(113, 12).Value = 2999
$ws.Cells.Item(113, 13).Value = 593.4614999999999
$ws.Cells.Item(113, 14).Value = -7339
$ws.Cells.Item(132, 8).Value = 2569.93
$ws.Cells.Item(132, 9).Value = 2362.0562
$ws.Cells.Item(132, 10).Value = 4251.8184
$ws.Cells.Item(132, 11).Value = 7086.1686
$ws.Cells.Item(132, 12).Value = 12755.4552
$ws.Cells.Item(132, 13).Value = -4556.1686
$ws.Cells.Item(132, 14).Value = -17815.4552
$ws.Cells.Item(136, 8).Value = 3200.111
$ws.Cells.Item(136, 9).Value = 2294.6978
$ws.Cells.Item(136, 11).Value = 6884.0934
$ws.Cells.Item(136, 13).Value = -4334.0934
$ws.Cells.Item(140, 8).Value = 55324.582
$ws.Cells.Item(140, 10).Value = 55324.582
$ws.Cells.Item(140, 12).Value = 55324.582
$ws.Cells.Item(140, 14).Value = -65684.58199999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 2758.9375
$ws.Cells.Item(62, 9).Value = 3144
$ws.Cells.Item(62, 10).Value = 2117.1667
$ws.Cells.Item(62, 11).Value = 3144
$ws.Cells.Item(62, 12).Value = 2117.1667
$ws.Cells.Item(62, 13).Value = -2520
$ws.Cells.Item(62, 14).Value = -3365.1667
$ws.Cells.Item(65, 8).Value = 2758.9375
$ws.Cells.Item(65, 9).Value = 3144
$ws.Cells.Item(65, 10).Value = 2117.1667
$ws.Cells.Item(65, 11).Value = 15720
$ws.Cells.Item(65, 12).Value = 10585.8335
$ws.Cells.Item(65, 13).Value = -12600
$ws.Cells.Item(65, 14).Value = -16825.8335
$ws.Cells.Item(81, 8).Value = 6826.231
$ws.Cells.Item(81, 9).Value = 7113
$ws.Cells.Item(81, 10).Value = 6491.6665
$ws.Cells.Item(81, 11).Value = 14226
$ws.Cells.Item(81, 12).Value = 12983.333
$ws.Cells.Item(81, 13).Value = -13165
$ws.Cells.Item(81, 14).Value = -15105.333
$ws.Cells.Item(84, 8).Value = 6826.231
$ws.Cells.Item(84, 9).Value = 7113
$ws.Cells.Item(84, 10).Value = 6491.6665
$ws.Cells.Item(84, 11).Value = 71130
$ws.Cells.Item(84, 12).Value = 64916.665
$ws.Cells.Item(84, 13).Value = -65826
$ws.Cells.Item(84, 14).Value = -75524.66500000001
$ws.Cells.Item(122, 8).Value = 3761.1785
$ws.Cells.Item(122, 9).Value = 2401.0588
$ws.Cells.Item(122, 10).Value = 5863.1816
$ws.Cells.Item(122, 11).Value = 7203.176399999999
$ws.Cells.Item(122, 12).Value = 17589.5448
$ws.Cells.Item(122, 13).Value = -4753.176399999999
$ws.Cells.Item(122, 14).Value = -22489.5448
$ws.Cells.Item(126, 8).Value = 2485.7144
$ws.Cells.Item(126, 9).Value = 2566.8333
$ws.Cells.Item(126, 11).Value = 7700.499899999999
$ws.Cells.Item(126, 13).Value = -5230.499899999999
$ws.Cells.Item(132, 8).Value = 1219.3422
$ws.Cells.Item(132, 9).Value = 917.4211
$ws.Cells.Item(132, 10).Value = 2125.1052
$ws.Cells.Item(132, 11).Value = 2752.2633
$ws.Cells.Item(132, 12).Value = 6375.3156
$ws.Cells.Item(132, 13).Value = -222.2633000000001
$ws.Cells.Item(132, 14).Value = -11435.3156
$ws.Cells.Item(136, 8).Value = 4825.98
$ws.Cells.Item(136, 10).Value = 9347.239
$ws.Cells.Item(136, 12).Value = 28041.717
$ws.Cells.Item(136, 14).Value = -33141.717
